$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 64

# Column A holds a date-shaped string ("2025/10/05"); Excel's usual
# autoconvert-to-date behavior would kick in on a plain .Value assignment,
# so force the cell to text first, then restore the default "Normal"
# style/number format once the literal text value has been stored.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2025/10/05"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "日"
$ws.Cells.Item($row, 3).Value = 16
$ws.Cells.Item($row, 4).Value = 55
